$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1267.8823
$ws.Range("I38").Value = 542.7778
$ws.Range("J38").Value = 2083.625
$ws.Range("K38").Value = 1628.3334
$ws.Range("L38").Value = 6250.875
$ws.Range("M38").Value = -1256.3334
$ws.Range("N38").Value = -6994.875

$ws.Range("H61").Value = 2355
$ws.Range("J61").Value = 11158.5
$ws.Range("L61").Value = 33475.5
$ws.Range("N61").Value = -33819.5

$ws.Range("H129").Value = 1209
$ws.Range("I129").Value = 901.4167
$ws.Range("J129").Value = 1282.82
$ws.Range("K129").Value = 2704.2501
$ws.Range("L129").Value = 3848.46
$ws.Range("M129").Value = 2295.7499
$ws.Range("N129").Value = -13848.46

$ws.Range("H137").Value = 2235.3333
$ws.Range("J137").Value = 1036.875
$ws.Range("L137").Value = 3110.625
$ws.Range("N137").Value = -8210.625

$ws.Range("H138").Value = 2084.597
$ws.Range("J138").Value = 2298.2856
$ws.Range("L138").Value = 6894.8568
$ws.Range("N138").Value = -17174.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 3000
$ws.Range("I25").Value = 3000
$ws.Range("K25").Value = 3000
$ws.Range("M25").Value = -2598

$ws.Range("H32").Value = 6552.71
$ws.Range("I32").Value = 5528.427
$ws.Range("J32").Value = 14840.091
$ws.Range("K32").Value = 5528.427
$ws.Range("L32").Value = 14840.091
$ws.Range("M32").Value = -5241.427
$ws.Range("N32").Value = -15414.091

$ws.Range("H35").Value = 837
$ws.Range("I35").Value = 837
$ws.Range("K35").Value = 837
$ws.Range("M35").Value = -431

$ws.Range("H74").Value = 5022.5557
$ws.Range("I74").Value = 6200.4614
$ws.Range("J74").Value = 1960
$ws.Range("K74").Value = 6200.4614
$ws.Range("L74").Value = 1960
$ws.Range("M74").Value = -5326.4614
$ws.Range("N74").Value = -3708

$ws.Range("H77").Value = 5022.5557
$ws.Range("I77").Value = 6200.4614
$ws.Range("J77").Value = 1960
$ws.Range("K77").Value = 31002.307
$ws.Range("L77").Value = 9800
$ws.Range("M77").Value = -26634.307
$ws.Range("N77").Value = -18536

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 863
$ws.Range("I37").Value = 863
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 863
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -726
$ws.Range("N37").ClearContents()

$ws.Range("H53").Value = 32500
$ws.Range("J53").Value = 32500
$ws.Range("L53").Value = 32500
$ws.Range("N53").Value = -33648

$ws.Range("H134").Value = 5494.8535
$ws.Range("I134").Value = 2760.4375
$ws.Range("J134").Value = 7244.88
$ws.Range("K134").Value = 8281.3125
$ws.Range("L134").Value = 21734.64
$ws.Range("M134").Value = -5746.3125
$ws.Range("N134").Value = -26804.64

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 9302.799999999999
$ws.Range("J43").Value = 9302.799999999999
$ws.Range("L43").Value = 9302.799999999999
$ws.Range("N43").Value = -9670.799999999999

$ws.Range("H58").Value = 6243.423
$ws.Range("I58").Value = 3496.2632
$ws.Range("J58").Value = 13700
$ws.Range("K58").Value = 3496.2632
$ws.Range("L58").Value = 13700
$ws.Range("M58").Value = -3293.2632
$ws.Range("N58").Value = -14106

$ws.Range("H101").Value = 9302.799999999999
$ws.Range("J101").Value = 9302.799999999999
$ws.Range("L101").Value = 9302.799999999999
$ws.Range("N101").Value = -15792.8

$ws.Range("H127").Value = 29890
$ws.Range("J127").Value = 39000
$ws.Range("L127").Value = 39000
$ws.Range("N127").Value = -48920

$ws.Range("H134").Value = 2588.3333
$ws.Range("I134").Value = 1200
$ws.Range("J134").Value = 2714.5454
$ws.Range("K134").Value = 3600
$ws.Range("L134").Value = 8143.6362
$ws.Range("M134").Value = -1065
$ws.Range("N134").Value = -13213.6362

$ws.Range("H136").Value = 6243.423
$ws.Range("I136").Value = 3496.2632
$ws.Range("J136").Value = 13700
$ws.Range("K136").Value = 10488.7896
$ws.Range("L136").Value = 41100
$ws.Range("M136").Value = -7938.7896
$ws.Range("N136").Value = -46200

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 574
$ws.Range("I107").Value = 350
$ws.Range("J107").Value = 798
$ws.Range("K107").Value = 1050
$ws.Range("L107").Value = 2394
$ws.Range("M107").Value = 870
$ws.Range("N107").Value = -6234

$ws.Range("H131").Value = 1046.234
$ws.Range("J131").Value = 1121.075
$ws.Range("L131").Value = 3363.225
$ws.Range("N131").Value = -13443.225

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 10785.5
$ws.Range("I99").Value = 7535.273
$ws.Range("K99").Value = 7535.273
$ws.Range("M99").Value = -5289.273

$ws.Range("H122").Value = 4670.3076
$ws.Range("I122").Value = 4539.25
$ws.Range("J122").Value = 4880
$ws.Range("K122").Value = 13617.75
$ws.Range("L122").Value = 14640
$ws.Range("M122").Value = -11167.75
$ws.Range("N122").Value = -19540

$ws.Range("H133").Value = 50170
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 50170
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 50170
$ws.Range("N133").Value = -60290
$ws.Range("M133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3547.5557
$ws.Range("I40").Value = 2969
$ws.Range("J40").Value = 5572.5
$ws.Range("K40").Value = 2969
$ws.Range("L40").Value = 5572.5
$ws.Range("M40").Value = -2833
$ws.Range("N40").Value = -5844.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 29700
$ws.Range("J105").Value = 29700
$ws.Range("L105").Value = 29700
$ws.Range("N105").Value = -36688

$ws.Range("H126").Value = 2306.162
$ws.Range("I126").Value = 2222.4285
$ws.Range("J126").Value = 2785.7273
$ws.Range("K126").Value = 6667.2855
$ws.Range("L126").Value = 8357.1819
$ws.Range("M126").Value = -4197.2855
$ws.Range("N126").Value = -13297.1819

$ws.Range("H136").Value = 928522.2
$ws.Range("I136").Value = 3706988.8
$ws.Range("J136").Value = 2366.6667
$ws.Range("K136").Value = 11120966.4
$ws.Range("L136").Value = 7100.000100000001
$ws.Range("M136").Value = -11118416.4
$ws.Range("N136").Value = -12200.0001
